$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert 5 new "FTT-IH-*" input sheets between "FTT-Fr" and "Time_Horizons"
#    Each sheet has the same 7-column header row and a single data row that
#    describes a fuel-tax (ISXn) variable.
# ---------------------------------------------------------------------------

$newSheetsInfo = @(
    @{ Name = "FTT-IH-CHI"; Code = "ISX1" },
    @{ Name = "FTT-IH-FBT"; Code = "ISX2" },
    @{ Name = "FTT-IH_MTM"; Code = "ISX3" },
    @{ Name = "FTT-IH-NMM"; Code = "ISX4" },
    @{ Name = "FTT-IH-OIS"; Code = "ISX5" }
)

$afterSheet = $wb.Worksheets.Item("FTT-Fr")

foreach ($info in $newSheetsInfo) {
    $ns = $wb.Worksheets.Add([System.Type]::Missing, $afterSheet)
    $ns.Name = $info.Name

    $ns.Range("A1").Value = "Variable name"
    $ns.Range("B1").Value = "Read in?"
    $ns.Range("C1").Value = "Code"
    $ns.Range("D1").Value = "Description"
    $ns.Range("E1").Value = "RowDim"
    $ns.Range("F1").Value = "ColDim"
    $ns.Range("G1").Value = "3DDim"

    $ns.Range("A2").Value = $info.Code
    $ns.Range("B2").Value = 1
    $ns.Range("C2").Value = 0
    $ns.Range("D2").Value = "FTT-Power historical generation"
    $ns.Range("E2").Value = "ITTI"
    $ns.Range("F2").Value = "TIME"
    $ns.Range("G2").Value = "RSHORTTI"

    $ns.Range("A3").Select()

    $afterSheet = $ns
}

# ---------------------------------------------------------------------------
# 2. Add the 5 new ISX1-5 rows (with their time horizon "tl_2000") at the
#    bottom of the "Time_Horizons" sheet.
# ---------------------------------------------------------------------------

$th = $wb.Worksheets.Item("Time_Horizons")

$th.Range("A61").Value = "ISX1"
$th.Range("B61").Value = "tl_2000"

$th.Range("A62").Value = "ISX2"
$th.Range("B62").Value = "tl_2000"

$th.Range("A63").Value = "ISX3"
$th.Range("B63").Value = "tl_2000"

$th.Range("A64").Value = "ISX4"
$th.Range("B64").Value = "tl_2000"

$th.Range("A65").Value = "ISX5"
$th.Range("B65").Value = "tl_2000"

# ---------------------------------------------------------------------------
# 3. Clear the (now unneeded) highlight fill on FTT-Fr!H4, H5 and F25.
# ---------------------------------------------------------------------------

$ftfr = $wb.Worksheets.Item("FTT-Fr")
$ftfr.Range("H4").Interior.Pattern = -4142
$ftfr.Range("H5").Interior.Pattern = -4142
$ftfr.Range("F25").Interior.Pattern = -4142

# ---------------------------------------------------------------------------
# 4. Time_Horizons becomes the active sheet / tab, matching the saved view.
# ---------------------------------------------------------------------------

$th.Activate()
$th.Range("J41").Select()
